$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new weekly price observation is inserted as row 142 (a new Monday of
# data for "Start Ruby" / "Primera"), which pushes all subsequent rows
# (old 142..171) down by one, ending with a new last row 172 that holds
# what used to be row 171's data. Excel's native row Insert() shifts the
# existing rows down and copies formatting (so column D keeps its date
# style) exactly like the author's edit.
$ws.Rows("142").Insert()

$ws.Range("A142").Value = 4
$ws.Range("B142").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C142").Value = "Los Lagos"
$ws.Range("D142").Value = 44511
$ws.Range("E142").Value = 10
$ws.Range("F142").Value = "Fruta"
$ws.Range("G142").Value = 100102
$ws.Range("H142").Value = "Cítricos"
$ws.Range("I142").Value = 100102006
$ws.Range("J142").Value = "Pomelo"
$ws.Range("K142").Value = "Start Ruby"
$ws.Range("L142").Value = "Primera"
$ws.Range("M142").Value = 60
$ws.Range("N142").Value = 11000
$ws.Range("O142").Value = 12000
$ws.Range("P142").Value = 11500
$ws.Range("Q142").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R142").Value = "Región de O'Higgins"
$ws.Range("S142").Value = 821
$ws.Range("T142").Value = 14
